# Add a new "2022-Q1" holdings sheet (positioned right before "总计"),
# populate it with the fund-holding data, and prepend a matching
# "2022-Q1" summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$templateSheet = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q1" sheet right before "总计"
# ---------------------------------------------------------------------
# NOTE: after Worksheets.Add(before), the handle passed as "before" gets
# reseated onto the newly inserted sheet, so "总计" must be re-resolved
# by name afterwards rather than reusing the old variable.
$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# Fund holding rows: code, name, fund size, stock position, position pct,
# holding value (billions CNY), position rank
$rows = @(
    @("003853", "金鹰信息产业股票A", "18.11", "89.63", "4.22", "0.7642", 5),
    @("515880", "国泰中证全指通信设备ETF", "16.24", "99.37", "3.38", "0.5489", 9),
    @("005885", "金鹰信息产业股票C", "6.45", "89.63", "4.22", "0.2722", 5),
    @("001809", "中信建投智信物联网灵活配置混合A", "5.79", "92.66", "3.34", "0.1934", 10),
    @("011685", "创金合信先进装备股票A", "0.73", "92.01", "9.84", "0.0718", 1),
    @("004636", "中信建投智信物联网灵活配置混合C", "1.66", "92.66", "3.34", "0.0554", 10),
    @("000714", "诺安稳健回报灵活配置混合A", "1.96", "64.55", "2.73", "0.0535", 3),
    @("002052", "诺安稳健回报灵活配置混合C", "1.63", "64.55", "2.73", "0.0445", 3),
    @("014133", "工银中证500六个月持有指数增强A", "3.07", "93.69", "0.92", "0.0282", 4),
    @("004351", "汇丰晋信珠三角区域发展混合", "0.51", "93.92", "4.54", "0.0232", 9),
    @("011686", "创金合信先进装备股票C", "0.17", "92.01", "9.84", "0.0167", 1),
    @("014134", "工银中证500六个月持有指数增强C", "1.12", "93.69", "0.92", "0.0103", 4),
    @("165524", "信诚中证智能家居指数（LOF）", "0.40", "93.89", "1.16", "0.0046", 9)
)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")

# --- values first (so the text-forcing NumberFormat trick only ever
#     touches cells that are still at the default style) ---
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2

    $newSheet.Cells.Item($excelRow, 1).Value = $r

    for ($c = 0; $c -lt 5; $c++) {
        $cell = $newSheet.Cells.Item($excelRow, $c + 2)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c]
        $cell.ClearFormats()
    }

    $valueCell = $newSheet.Cells.Item($excelRow, 7)
    $valueCell.NumberFormat = "@"
    $valueCell.Value = $row[5]
    $valueCell.ClearFormats()

    $newSheet.Cells.Item($excelRow, 8).Value = $row[6]
}

# --- now stamp the styles used by the other quarterly sheets onto the
#     header row and the index column, without touching the values ---
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2:A4").Copy()
$newSheet.Range("A2:A14").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Insert the new "2022-Q1" row at the top of the "总计" sheet
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 13
$totalSheet.Range("D2").Value = 2.09

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
